# 16.3.1 workbook update: add a 2023 data column (E) alongside the
# existing 2018 column (D), refresh the footnote to mention the new
# 2018/2023 cluster survey, and bold the "Urbanisation" header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# 1. Footnote row (row 9): append the 2023 survey reference next to
#    the existing 2018 reference, in all three languages.
# ---------------------------------------------------------------
$ws.Range("A9").Value = " Көп көрсөткүчтүү кластердик изилдөөнүн маалыматтары боюнча, 2018-ж., 2023-ж."
$ws.Range("B9").Value = "По данным кластерного обследования по многим показателям, 2018г., 2023г."
$ws.Range("C9").Value = "According to the cluster survey in many respects, 2018, 2023."

# ---------------------------------------------------------------
# 2. New column E holding the 2023 figures, mirroring column D's
#    layout/formatting (year header, overall value, blank section
#    header, and two "-" placeholders for the urban/rural split).
# ---------------------------------------------------------------
$ws.Range("D4").Copy() | Out-Null
$ws.Range("E4").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("E4").Value = 2023

$ws.Range("D5").Copy() | Out-Null
$ws.Range("E5").PasteSpecial(-4122) | Out-Null
$ws.Range("E5").Value = 38

$ws.Range("D6").Copy() | Out-Null
$ws.Range("E6").PasteSpecial(-4122) | Out-Null

$ws.Range("D7").Copy() | Out-Null
$ws.Range("E7").PasteSpecial(-4122) | Out-Null
$ws.Range("E7").Value = "-"
$ws.Range("E7").HorizontalAlignment = -4152      # xlRight

$ws.Range("D8").Copy() | Out-Null
$ws.Range("E8").PasteSpecial(-4122) | Out-Null
$ws.Range("E8").Value = "-"
$ws.Range("E8").HorizontalAlignment = -4152      # xlRight

$ws.Application.CutCopyMode = $false

# ---------------------------------------------------------------
# 3. Bold the "Urbanisation" section-header row (now spanning
#    columns A:E after the new column was added).
# ---------------------------------------------------------------
$ws.Range("A6:E6").Font.Bold = $true
